$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2159311196958242
$ws.Range("D2").Value = 0.1740869811143284
$ws.Range("E2").Value = 0.1584551540436152
$ws.Range("F2").Value = 1.586058537684849
$ws.Range("G2").Value = 0.002415559989362569
$ws.Range("I2").Value = 0.3757315010515043
$ws.Range("J2").Value = 0.1913305446093858
$ws.Range("N2").Value = 1.538965158849294
$ws.Range("O2").Value = 3.975293305671926
$ws.Range("B3").Value = 0.1913949238933697
$ws.Range("D3").Value = 0.1738283514950183
$ws.Range("E3").Value = 0.1560976871575157
$ws.Range("F3").Value = 1.537274724656484
$ws.Range("G3").Value = 0.002420193372685223
$ws.Range("I3").Value = 0.3656238581545637
$ws.Range("J3").Value = 0.1859881629284246
$ws.Range("N3").Value = 1.440536212657179
$ws.Range("O3").Value = 3.829925653828582
$ws.Range("B4").Value = 0.1762869034279646
$ws.Range("D4").Value = 0.1737352727050379
$ws.Range("E4").Value = 0.1547365595057464
$ws.Range("F4").Value = 1.508292599236057
$ws.Range("G4").Value = 0.002423190390540851
$ws.Range("I4").Value = 0.3596127975951688
$ws.Range("J4").Value = 0.1828259976756428
$ws.Range("N4").Value = 1.380350891855613
$ws.Range("O4").Value = 3.743161388050623
$ws.Range("B5").Value = 0.1701201473106835
$ws.Range("D5").Value = 0.1737138821688973
$ws.Range("E5").Value = 0.1542035670625133
$ws.Range("F5").Value = 1.496725533069323
$ws.Range("G5").Value = 0.002424450076845931
$ws.Range("I5").Value = 0.3572128864359954
$ws.Range("J5").Value = 0.1815669644110329
$ws.Range("N5").Value = 1.355891405433653
$ws.Range("O5").Value = 3.708428181857073
$ws.Range("B6").Value = 0.1690955695415823
$ws.Range("D6").Value = 0.1737113296674551
$ws.Range("E6").Value = 0.1541163722461825
$ws.Range("F6").Value = 1.494819505439395
$ws.Range("G6").Value = 0.002424661568431892
$ws.Range("I6").Value = 0.3568173988915788
$ws.Range("J6").Value = 0.1813596857490793
$ws.Range("N6").Value = 1.351834048475041
$ws.Range("O6").Value = 3.702698365185597
$ws.Range("B7").Value = 0.1762037765295901
$ws.Range("D7").Value = 0.1737349172359117
$ws.Range("E7").Value = 0.1547292836604051
$ws.Range("F7").Value = 1.508135617428664
$ws.Range("G7").Value = 0.002423207223543514
$ws.Range("I7").Value = 0.3595802297560198
$ws.Range("J7").Value = 0.1828088983185268
$ws.Range("N7").Value = 1.38002074861987
$ws.Range("O7").Value = 3.742690441195975
$ws.Range("B8").Value = 0.2074802658001431
$ws.Range("D8").Value = 0.1739841729763967
$ws.Range("E8").Value = 0.1576243489920373
$ws.Range("F8").Value = 1.569035736272298
$ws.Range("G8").Value = 0.002417126089359947
$ws.Range("I8").Value = 0.3722062484302811
$ws.Range("J8").Value = 0.1894639112100833
$ws.Range("N8").Value = 1.504976964632107
$ws.Range("O8").Value = 3.92465161631327
$ws.Range("B9").Value = 0.2684501888315936
$ws.Range("D9").Value = 0.1749940350993739
$ws.Range("E9").Value = 0.1639891454149733
$ws.Range("F9").Value = 1.696215314610981
$ws.Range("G9").Value = 0.002406401938415241
$ws.Range("I9").Value = 0.3984893586955494
$ws.Range("J9").Value = 0.2034575054670285
$ws.Range("N9").Value = 1.751860150726543
$ws.Range("O9").Value = 4.301396925444635
$ws.Range("B10").Value = 0.3129948669844396
$ws.Range("D10").Value = 0.1760534290396052
$ws.Range("E10").Value = 0.169088474725001
$ws.Range("F10").Value = 1.794459758914059
$ws.Range("G10").Value = 0.002399246738317698
$ws.Range("I10").Value = 0.4186982167615625
$ws.Range("J10").Value = 0.2143235631693869
$ws.Range("N10").Value = 1.93419659562025
$ws.Range("O10").Value = 4.590573542432821
$ws.Range("B11").Value = 0.33319936175468
$ws.Range("D11").Value = 0.1766042903173286
$ws.Range("E11").Value = 0.1715010782485393
$ws.Range("F11").Value = 1.840215008519493
$ws.Range("G11").Value = 0.002396147058881212
$ws.Range("I11").Value = 0.4280806917633839
$ws.Range("J11").Value = 0.2193960977187146
$ws.Range("N11").Value = 2.017317443662478
$ws.Range("O11").Value = 4.724870106639514
$ws.Range("B12").Value = 0.3408412057247006
$ws.Range("D12").Value = 0.1768227941903859
$ws.Range("E12").Value = 0.1724280870050521
$ws.Range("F12").Value = 1.857695457956879
$ws.Range("G12").Value = 0.00239499548138893
$ws.Range("I12").Value = 0.4316602591565939
$ws.Range("J12").Value = 0.221335716487971
$ws.Range("N12").Value = 2.048815004509549
$ws.Range("O12").Value = 4.776123705914472
$ws.Range("B13").Value = 0.3391958157679937
$ws.Range("D13").Value = 0.1767752950894987
$ws.Range("E13").Value = 0.1722278419399359
$ws.Range("F13").Value = 1.853923867927108
$ws.Range("G13").Value = 0.002395242508514029
$ws.Range("I13").Value = 0.4308881605543888
$ws.Range("J13").Value = 0.2209171483289651
$ws.Range("N13").Value = 2.042030543333738
$ws.Range("O13").Value = 4.765067556405256
$ws.Range("B14").Value = 0.3338282482071406
$ws.Range("D14").Value = 0.1766220683052069
$ws.Range("E14").Value = 0.1715770748056329
$ws.Range("F14").Value = 1.841650044738657
$ws.Range("G14").Value = 0.00239605187356056
$ws.Range("I14").Value = 0.4283746547068787
$ws.Range("J14").Value = 0.2195552945535866
$ws.Range("N14").Value = 2.019908357050213
$ws.Range("O14").Value = 4.729078768254283
$ws.Range("B15").Value = 0.3305392456628908
$ws.Range("D15").Value = 0.1765295021117765
$ws.Range("E15").Value = 0.1711802087939134
$ws.Range("F15").Value = 1.834152050347342
$ws.Range("G15").Value = 0.002396550521006046
$ws.Range("I15").Value = 0.4268385096919047
$ws.Range("J15").Value = 0.2187235669912297
$ws.Range("N15").Value = 2.006360570313291
$ws.Range("O15").Value = 4.707086574762002
$ws.Range("B16").Value = 0.3116732098856971
$ws.Range("D16").Value = 0.1760188154244986
$ws.Range("E16").Value = 0.1689326789279946
$ws.Range("F16").Value = 1.791491030024417
$ws.Range("G16").Value = 0.00239945242344497
$ws.Range("I16").Value = 0.4180888088535895
$ws.Range("J16").Value = 0.2139946788174569
$ws.Range("N16").Value = 1.928767661017872
$ws.Range("O16").Value = 4.581852524519547
$ws.Range("B17").Value = 0.3000838993125114
$ws.Range("D17").Value = 0.1757231760749818
$ws.Range("E17").Value = 0.167577719709449
$ws.Range("F17").Value = 1.765592947740373
$ws.Range("G17").Value = 0.002401272324867566
$ws.Range("I17").Value = 0.4127692216154415
$ws.Range("J17").Value = 0.2111269263642015
$ws.Range("N17").Value = 1.88120914429345
$ws.Range("O17").Value = 4.505731717013816
$ws.Range("B18").Value = 0.2934125106850161
$ws.Range("D18").Value = 0.1755596213569888
$ws.Range("E18").Value = 0.1668071234688995
$ws.Range("F18").Value = 1.750797082196556
$ws.Range("G18").Value = 0.002402333704617743
$ws.Range("I18").Value = 0.4097273953840173
$ws.Range("J18").Value = 0.2094896528478927
$ws.Range("N18").Value = 1.85387144156806
$ws.Range("O18").Value = 4.462207568522672
$ws.Range("B19").Value = 0.2911527646049308
$ws.Range("D19").Value = 0.1755053593357232
$ws.Range("E19").Value = 0.166547712491564
$ws.Range("F19").Value = 1.745804605181775
$ws.Range("G19").Value = 0.00240269558458022
$ws.Range("I19").Value = 0.4087005711172083
$ws.Range("J19").Value = 0.2089373882823509
$ws.Range("N19").Value = 1.844618345498787
$ws.Range("O19").Value = 4.447515346699447
$ws.Range("B20").Value = 0.3013181776851468
$ws.Range("D20").Value = 0.1757539758419355
$ws.Range("E20").Value = 0.167721052465474
$ws.Range("F20").Value = 1.768339483154278
$ws.Range("G20").Value = 0.002401077080864969
$ws.Range("I20").Value = 0.4133336567293782
$ws.Range("J20").Value = 0.2114309417489437
$ws.Range("N20").Value = 1.886270128093656
$ws.Range("O20").Value = 4.513808127811615
$ws.Range("B21").Value = 0.3354050867553156
$ws.Range("D21").Value = 0.1766668059812773
$ws.Range("E21").Value = 0.1717678565883389
$ws.Range("F21").Value = 1.845250978317722
$ws.Range("G21").Value = 0.002395813541770542
$ws.Range("I21").Value = 0.429112214418673
$ws.Range("J21").Value = 0.2199547937791095
$ws.Range("N21").Value = 2.026405625114307
$ws.Range("O21").Value = 4.739638709769963
$ws.Range("B22").Value = 0.3576291578720543
$ws.Range("D22").Value = 0.1773211216162025
$ws.Range("E22").Value = 0.1744908541518555
$ws.Range("F22").Value = 1.896414790630246
$ws.Range("G22").Value = 0.002392502880666877
$ws.Range("I22").Value = 0.4395793436413413
$ws.Range("J22").Value = 0.225635036445496
$ws.Range("N22").Value = 2.118115898853375
$ws.Range("O22").Value = 4.889555756343043
$ws.Range("B23").Value = 0.3457728835225851
$ws.Range("D23").Value = 0.1769666215723475
$ws.Range("E23").Value = 0.1730303690042234
$ws.Range("F23").Value = 1.869025219674057
$ws.Range("G23").Value = 0.002394258045721026
$ws.Range("I23").Value = 0.4339788619098002
$ws.Range("J23").Value = 0.2225933271110279
$ws.Range("N23").Value = 2.069158297353965
$ws.Range("O23").Value = 4.809328534290898
$ws.Range("B24").Value = 0.3007601869820462
$ws.Range("D24").Value = 0.1757400312936568
$ws.Range("E24").Value = 0.1676562256042899
$ws.Range("F24").Value = 1.76709748398801
$ws.Range("G24").Value = 0.002401165303691459
$ws.Range("I24").Value = 0.4130784242511751
$ws.Range("J24").Value = 0.2112934607855692
$ws.Range("N24").Value = 1.883982043575344
$ws.Range("O24").Value = 4.510156039135268
$ws.Range("B25").Value = 0.2519980490933449
$ws.Range("D25").Value = 0.1746650798694702
$ws.Range("E25").Value = 0.1621932848350376
$ws.Range("F25").Value = 1.660971659900625
$ws.Range("G25").Value = 0.002409175396314761
$ws.Range("I25").Value = 0.3912191043998661
$ws.Range("J25").Value = 0.1995698711105831
$ws.Range("N25").Value = 1.684892220676289
$ws.Range("O25").Value = 4.197320322167002
